$d = $word.ActiveDocument

# --- Body: FECHA date "27/mayo/2022" -> "1/junio/2022" ---
$body = $d.Content
$body.Find.ClearFormatting()
$body.Find.Replacement.ClearFormatting()
$body.Find.Execute("27", $true, $true, $false, $false, $false, $true, 1, $false, "1", 2)

$body2 = $d.Content
$body2.Find.ClearFormatting()
$body2.Find.Replacement.ClearFormatting()
$body2.Find.Execute("mayo", $true, $true, $false, $false, $false, $true, 1, $false, "junio", 2)

# --- Header: PERIODO "2022-05-26  -  2022-05-27" -> "2022-06-01  -  2022-06-02" ---
$section = $d.Sections.Item(1)
$header = $section.Headers.Item(1)

$headerRange1 = $header.Range
$headerRange1.Find.ClearFormatting()
$headerRange1.Find.Replacement.ClearFormatting()
$headerRange1.Find.Execute("2022-05-26", $true, $false, $false, $false, $false, $true, 1, $false, "2022-06-01", 2)

$headerRange2 = $header.Range
$headerRange2.Find.ClearFormatting()
$headerRange2.Find.Replacement.ClearFormatting()
$headerRange2.Find.Execute("2022-05-27", $true, $false, $false, $false, $false, $true, 1, $false, "2022-06-02", 2)
